# start of season update
# Correct the slightly-rounded timestamp on the last existing log row and
# append the new season's log entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130: refresh the timestamp (tiny precision correction).
$ws.Cells.Item(130, 1).Value = 45443.33009739583

# Row 131: new log entry.
$ws.Cells.Item(131, 1).Value = 45524.55052809028
$ws.Cells.Item(131, 2).Value = 5
$ws.Cells.Item(131, 3).Value = 130
$ws.Cells.Item(131, 4).Value = "data extract routine run"

# Row 132: new log entry.
$ws.Cells.Item(132, 1).Value = 45525.31197869213
$ws.Cells.Item(132, 2).Value = 4
$ws.Cells.Item(132, 3).Value = 131
$ws.Cells.Item(132, 4).Value = "data extract routine run"

# Row 133: new log entry.
$ws.Cells.Item(133, 1).Value = 45525.33312254546
$ws.Cells.Item(133, 2).Value = 7
$ws.Cells.Item(133, 3).Value = 132
$ws.Cells.Item(133, 4).Value = "data upload routine run"
